$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: re-curated dimension/measure labels
$ws.Range("A2").Value = "sdmx-dimension:refArea"
$ws.Range("B2").Value = "iaest-measure:nivel-estudios"

# Row 3: role labels swap (municipio column becomes the dimension, nivel-estudios becomes the measure)
$ws.Range("A3").Value = "dim"
$ws.Range("B3").Value = "medida"

# Row 4: datatype / URI columns
$ws.Range("A4").Value = "URI-Municipio"
$ws.Range("B4").Value = "xsd:int"

# Row 5 no longer applies (mapping file reference removed) - delete entire row
$ws.Rows(5).Delete()
